# Update countries & provincias Spain
# Applies the 2020-05-17 21:05 data refresh to the "Pais" sheet:
#  - updates the "last updated" timestamp
#  - updates totals for Estados Unidos, Jordania, Republica del Chad
#  - refreshes Sudan del Sur / Yemen / Gambia figures and re-sorts them
#    ahead of Ruanda / Bermudas / Burundi respectively (ranking swaps
#    caused by the new totals)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header timestamp
$ws.Range("A1").Value = "Datos actualizados a 17 de Mayo de 2020 a las 21:05"

# Estados Unidos (row 4) - updated totals
$ws.Cells.Item(4, 2).Value = 1516004
$ws.Cells.Item(4, 3).Value = 8231
$ws.Cells.Item(4, 5).Value = 1084878
$ws.Cells.Item(4, 7).Value = 286
$ws.Cells.Item(4, 8).Value = 90399

# Jordania (row 122) - updated totals
$ws.Cells.Item(122, 2).Value = 613
$ws.Cells.Item(122, 3).Value = 6
$ws.Cells.Item(122, 4).Value = 408
$ws.Cells.Item(122, 5).Value = 196

# Republica del Chad (row 129) - updated totals
$ws.Cells.Item(129, 2).Value = 503
$ws.Cells.Item(129, 3).Value = 29
$ws.Cells.Item(129, 4).Value = 117
$ws.Cells.Item(129, 5).Value = 333
$ws.Cells.Item(129, 7).Value = 3
$ws.Cells.Item(129, 8).Value = 53

# Row 146/147: Sudan del Sur overtakes Ruanda with fresh numbers
$ws.Cells.Item(146, 1).Value = "Sudan del Sur"
$ws.Cells.Item(146, 2).Value = 290
$ws.Cells.Item(146, 3).Value = 54
$ws.Cells.Item(146, 4).Value = 4
$ws.Cells.Item(146, 5).Value = 282
$ws.Cells.Item(146, 6).Value = 0
$ws.Cells.Item(146, 7).Value = 0
$ws.Cells.Item(146, 8).Value = 4

$ws.Cells.Item(147, 1).Value = "Ruanda"
$ws.Cells.Item(147, 2).Value = 289
$ws.Cells.Item(147, 3).Value = 0
$ws.Cells.Item(147, 4).Value = 178
$ws.Cells.Item(147, 5).Value = 111
$ws.Cells.Item(147, 6).Value = 0
$ws.Cells.Item(147, 7).Value = 0
$ws.Cells.Item(147, 8).Value = 0

# Row 161/162: Yemen overtakes Bermudas with fresh numbers
$ws.Cells.Item(161, 1).Value = "Yemen"
$ws.Cells.Item(161, 2).Value = 128
$ws.Cells.Item(161, 3).Value = 6
$ws.Cells.Item(161, 4).Value = 1
$ws.Cells.Item(161, 5).Value = 107
$ws.Cells.Item(161, 6).Value = 0
$ws.Cells.Item(161, 7).Value = 2
$ws.Cells.Item(161, 8).Value = 20

$ws.Cells.Item(162, 1).Value = "Bermudas"
$ws.Cells.Item(162, 2).Value = 123
$ws.Cells.Item(162, 3).Value = 0
$ws.Cells.Item(162, 4).Value = 73
$ws.Cells.Item(162, 5).Value = 41
$ws.Cells.Item(162, 6).Value = 0
$ws.Cells.Item(162, 7).Value = 0
$ws.Cells.Item(162, 8).Value = 9

# Row 189/190: Gambia overtakes Burundi with fresh numbers
$ws.Cells.Item(189, 1).Value = "Gambia"
$ws.Cells.Item(189, 2).Value = 23
$ws.Cells.Item(189, 3).Value = 0
$ws.Cells.Item(189, 4).Value = 12
$ws.Cells.Item(189, 5).Value = 10
$ws.Cells.Item(189, 6).Value = 0
$ws.Cells.Item(189, 7).Value = 0
$ws.Cells.Item(189, 8).Value = 1

$ws.Cells.Item(190, 1).Value = "Burundi"
$ws.Cells.Item(190, 2).Value = 23
$ws.Cells.Item(190, 3).Value = 8
$ws.Cells.Item(190, 4).Value = 15
$ws.Cells.Item(190, 5).Value = 7
$ws.Cells.Item(190, 6).Value = 0
$ws.Cells.Item(190, 7).Value = 0
$ws.Cells.Item(190, 8).Value = 1
